# "Generate Report for Handback"
#
# - "Ready for handoff" status becomes "Handback transform failed" wherever
#   it appears (Overview!E3/F3, zh-cn!C3, de-de!C3).
# - The (previously empty) "Error Detail" column (P) on the zh-cn and de-de
#   report sheets gets a diagnostic message for row 3, and that column is
#   widened to fit (40 characters, matching the other long-text columns).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = "Handback file name: wbopotxm.e3e is different with handoff file name: 40f37778-5c48-4421-b71c-3a87001b5be1.bb54d3ac3d0455dc5e7ee8a2ba6e441b9ffcdbd4.zh-cn."
# Excel's ColumnWidth (characters) round-trips through pixel units: the
# OOXML <col width> ends up ~0.8333 (5/6) wider than the value assigned
# here, so back off by that amount to land exactly on width="40".
$zhcn.Range("P3").ColumnWidth = 40 - 5/6

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = "Handback file name: wbopotxm.e3e is different with handoff file name: 40f37778-5c48-4421-b71c-3a87001b5be1.bb54d3ac3d0455dc5e7ee8a2ba6e441b9ffcdbd4.de-de."
$dede.Range("P3").ColumnWidth = 40 - 5/6
